$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- Sheet ALC ---
# Row 5
$ws1.Range("H5").Value = 133.375
$ws1.Range("I5").Value = 141.16667
$ws1.Range("J5").Value = 110
$ws1.Range("K5").Value = 141.16667
$ws1.Range("L5").Value = 110
$ws1.Range("M5").Value = -26.16667000000001
$ws1.Range("N5").Value = -340
# Row 31
$ws1.Range("H31").Value = 1900
$ws1.Range("I31").Value = 1000
$ws1.Range("J31").Value = 2800
$ws1.Range("K31").Value = 3000
$ws1.Range("L31").Value = 8400
$ws1.Range("M31").Value = -2770
$ws1.Range("N31").Value = -8860
# Row 38
$ws1.Range("H38").Value = 350.25

# --- Sheet ARM ---
# Row 30
$ws2.Range("H30").Value = 50009.668
$ws2.Range("I30").Value = 10009
$ws2.Range("J30").Value = 70010
$ws2.Range("K30").Value = 10009
$ws2.Range("L30").Value = 70010
$ws2.Range("M30").Value = -9859
$ws2.Range("N30").Value = -70310
# Row 32
$ws2.Range("H32").Value = 330216.75
$ws2.Range("I32").Value = 393999.62
$ws2.Range("J32").Value = 11302.267
$ws2.Range("K32").Value = 393999.62
$ws2.Range("L32").Value = 11302.267
$ws2.Range("M32").Value = -393712.62
$ws2.Range("N32").Value = -11876.267
# Row 43
$ws2.Range("H43").Value = 10835.25
$ws2.Range("I43").Value = 9170.5
$ws2.Range("J43").Value = 12500
$ws2.Range("K43").Value = 9170.5
$ws2.Range("L43").Value = 12500
$ws2.Range("M43").Value = -8857.5
$ws2.Range("N43").Value = -13126
# Row 45
$ws2.Range("H45").Value = 3421.5293
$ws2.Range("I45").Value = 2590.75
$ws2.Range("J45").Value = 4160
$ws2.Range("K45").Value = 2590.75
$ws2.Range("L45").Value = 4160
$ws2.Range("M45").Value = -2213.75
$ws2.Range("N45").Value = -4914
# Row 61
$ws2.Range("H61").Value = 7577989
$ws2.Range("I61").Value = 19609226
$ws2.Range("K61").Value = 19609226
$ws2.Range("M61").Value = -19609014
# Row 74
$ws2.Range("H74").Value = 1602.4546
$ws2.Range("I74").Value = 1341
$ws2.Range("J74").Value = 1700.5
$ws2.Range("K74").Value = 1341
$ws2.Range("L74").Value = 1700.5
$ws2.Range("M74").Value = -467
$ws2.Range("N74").Value = -3448.5
# Row 77
$ws2.Range("H77").Value = 1602.4546
$ws2.Range("I77").Value = 1341
$ws2.Range("J77").Value = 1700.5
$ws2.Range("K77").Value = 6705
$ws2.Range("L77").Value = 8502.5
$ws2.Range("M77").Value = -2337
$ws2.Range("N77").Value = -17238.5
# Row 122
$ws2.Range("H122").Value = 1533.3334
$ws2.Range("I122").Value = 1533.3334
$ws2.Range("J122").Value = 0
$ws2.Range("K122").Value = 4600.0002
$ws2.Range("L122").Value = 0
$ws2.Range("M122").Value = -2150.0002
$ws2.Range("N122").ClearContents() | Out-Null
# Row 136
$ws2.Range("H136").Value = 7577989
$ws2.Range("I136").Value = 19609226
$ws2.Range("K136").Value = 58827678
$ws2.Range("M136").Value = -58825128

# --- Sheet BSM ---
# Row 108
$ws3.Range("H108").Value = 98684
$ws3.Range("J108").Value = 98684
$ws3.Range("L108").Value = 98684
$ws3.Range("N108").Value = -106364

# --- Sheet CRP ---
# Row 19
$ws4.Range("H19").Value = 3558.7646
$ws4.Range("I19").Value = 499.93332
$ws4.Range("J19").Value = 26500
$ws4.Range("K19").Value = 499.93332
$ws4.Range("L19").Value = 26500
$ws4.Range("M19").Value = -329.93332
$ws4.Range("N19").Value = -26840
# Row 24
$ws4.Range("H24").Value = 3558.7646
$ws4.Range("I24").Value = 499.93332
$ws4.Range("J24").Value = 26500
$ws4.Range("K24").Value = 499.93332
$ws4.Range("L24").Value = 26500
$ws4.Range("M24").Value = -329.93332
$ws4.Range("N24").Value = -26840
# Row 31
$ws4.Range("H31").Value = 5858.514
$ws4.Range("I31").Value = 1663.75
$ws4.Range("J31").Value = 8047.087
$ws4.Range("K31").Value = 1663.75
$ws4.Range("L31").Value = 8047.087
$ws4.Range("M31").Value = -1368.75
$ws4.Range("N31").Value = -8637.087
# Row 34
$ws4.Range("H34").Value = 5858.514
$ws4.Range("I34").Value = 1663.75
$ws4.Range("J34").Value = 8047.087
$ws4.Range("K34").Value = 1663.75
$ws4.Range("L34").Value = 8047.087
$ws4.Range("M34").Value = -1461.75
$ws4.Range("N34").Value = -8451.087
# Row 59
$ws4.Range("H59").Value = 27830.834
$ws4.Range("J59").Value = 27830.834
$ws4.Range("L59").Value = 27830.834
$ws4.Range("N59").Value = -30120.834
# Row 60
$ws4.Range("H60").Value = 12513.333
$ws4.Range("I60").Value = 9093
$ws4.Range("K60").Value = 9093
$ws4.Range("M60").Value = -8582
# Row 68
$ws4.Range("H68").Value = 28199.4
$ws4.Range("J68").Value = 28199.4
$ws4.Range("L68").Value = 28199.4
$ws4.Range("N68").Value = -29697.4
# Row 71
$ws4.Range("H71").Value = 28199.4
$ws4.Range("J71").Value = 28199.4
$ws4.Range("L71").Value = 84598.20000000001
$ws4.Range("N71").Value = -92086.20000000001
# Row 134
$ws4.Range("H134").Value = 1681.4445
$ws4.Range("I134").Value = 1022.1667
$ws4.Range("J134").Value = 3000
$ws4.Range("K134").Value = 3066.5001
$ws4.Range("L134").Value = 9000
$ws4.Range("M134").Value = -531.5001000000002
$ws4.Range("N134").Value = -14070

# --- Sheet CUL ---
# Row 34
$ws5.Range("H34").Value = 13889331
$ws5.Range("J34").Value = 16667149
$ws5.Range("L34").Value = 50001447
$ws5.Range("N34").Value = -50001615
# Row 50
$ws5.Range("H50").Value = 19608256
$ws5.Range("I50").Value = 418.18182
$ws5.Range("K50").Value = 1254.54546
$ws5.Range("M50").Value = -773.54546
# Row 53
$ws5.Range("H53").Value = 19608256
$ws5.Range("I53").Value = 418.18182
$ws5.Range("K53").Value = 1254.54546
$ws5.Range("M53").Value = -773.54546
# Row 55
$ws5.Range("H55").Value = 2074.0588
$ws5.Range("J55").Value = 2074.0588
$ws5.Range("L55").Value = 6222.176399999999
$ws5.Range("N55").Value = -6576.176399999999
# Row 68
$ws5.Range("H68").Value = 1402.8977
$ws5.Range("J68").Value = 1628.9231
$ws5.Range("L68").Value = 4886.7693
$ws5.Range("N68").Value = -6508.7693
# Row 71
$ws5.Range("H71").Value = 1402.8977
$ws5.Range("J71").Value = 1628.9231
$ws5.Range("L71").Value = 14660.3079
$ws5.Range("N71").Value = -22772.3079
# Row 113
$ws5.Range("H113").Value = 864.5
$ws5.Range("I113").Value = 513.48
$ws5.Range("J113").Value = 1380.7059
$ws5.Range("K113").Value = 1540.44
$ws5.Range("L113").Value = 4142.1177
$ws5.Range("M113").Value = 629.5599999999999
$ws5.Range("N113").Value = -8482.117699999999

# --- Sheet GSM ---
# Row 122
$ws6.Range("H122").Value = 3716.04
$ws6.Range("I122").Value = 2851.0833
$ws6.Range("K122").Value = 8553.249899999999
$ws6.Range("M122").Value = -6103.249899999999

# --- Sheet LTW ---
# Row 7
$ws7.Range("H7").Value = 43480630
$ws7.Range("I7").Value = 62502212
$ws7.Range("J7").Value = 2729.2856
$ws7.Range("K7").Value = 62502212
$ws7.Range("L7").Value = 2729.2856
$ws7.Range("M7").Value = -62502100
$ws7.Range("N7").Value = -2953.2856
# Row 40
$ws7.Range("H40").Value = 45457204
$ws7.Range("I40").Value = 55557664
$ws7.Range("J40").Value = 5125
$ws7.Range("K40").Value = 55557664
$ws7.Range("L40").Value = 5125
$ws7.Range("M40").Value = -55557528
$ws7.Range("N40").Value = -5397
# Row 122
$ws7.Range("H122").Value = 3502.6758
$ws7.Range("I122").Value = 2318.1365
$ws7.Range("J122").Value = 5240
$ws7.Range("K122").Value = 6954.4095
$ws7.Range("L122").Value = 15720
$ws7.Range("M122").Value = -4504.4095
$ws7.Range("N122").Value = -20620
# Row 126
$ws7.Range("H126").Value = 43480630
$ws7.Range("I126").Value = 62502212
$ws7.Range("J126").Value = 2729.2856
$ws7.Range("K126").Value = 187506636
$ws7.Range("L126").Value = 8187.8568
$ws7.Range("M126").Value = -187504166
$ws7.Range("N126").Value = -13127.8568

# --- Sheet WVR ---
# Row 122
$ws8.Range("H122").Value = 1789.95
$ws8.Range("I122").Value = 1471.3572
$ws8.Range("J122").Value = 2533.3333
$ws8.Range("K122").Value = 4414.071599999999
$ws8.Range("L122").Value = 7599.999899999999
$ws8.Range("M122").Value = -1964.071599999999
$ws8.Range("N122").Value = -12499.9999
# Row 136
$ws8.Range("H136").Value = 2871.186
$ws8.Range("I136").Value = 2780.64
$ws8.Range("J136").Value = 2996.9443
$ws8.Range("K136").Value = 8341.92
$ws8.Range("L136").Value = 8990.832900000001
$ws8.Range("M136").Value = -5791.92
$ws8.Range("N136").Value = -14090.8329
